# Apply updated TPM-derived statistics to the Efnb2-Pecam1 LR-pairs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 45.76217133333333
$ws.Range("H2").Value = 137.286514
$ws.Range("I2").Value = 0.6763939203605134
$ws.Range("J2").Value = 0.6763939203605135
$ws.Range("M2").Value = 361.2779286666667
$ws.Range("N2").Value = 1083.833786
$ws.Range("O2").Value = 0.9679392703861037
$ws.Range("P2").Value = 0.9679392703861038
$ws.Range("Q2").Value = 16532.86247059578
$ws.Range("R2").Value = 148795.762235362
$ws.Range("S2").Value = 0.6547082377673517
$ws.Range("T2").Value = 0.6547082377673519

# Row 3
$ws.Range("G3").Value = 45.76217133333333
$ws.Range("H3").Value = 137.286514
$ws.Range("I3").Value = 0.6763939203605134
$ws.Range("J3").Value = 0.6763939203605135
$ws.Range("O3").Value = 0.015995373883918
$ws.Range("P3").Value = 0.015995373883918
$ws.Range("Q3").Value = 273.2085830995254
$ws.Range("R3").Value = 2458.877247895728
$ws.Range("S3").Value = 0.01081917364897547
$ws.Range("T3").Value = 0.01081917364897547

# Row 4
$ws.Range("G4").Value = 45.76217133333333
$ws.Range("H4").Value = 137.286514
$ws.Range("I4").Value = 0.6763939203605134
$ws.Range("J4").Value = 0.6763939203605135
$ws.Range("M4").Value = 3.890485666666667
$ws.Range("N4").Value = 11.671457
$ws.Range("O4").Value = 0.01042342628440887
$ws.Range("P4").Value = 0.01042342628440887
$ws.Range("Q4").Value = 178.0370716478776
$ws.Range("R4").Value = 1602.333644830898
$ws.Range("S4").Value = 0.007050342168100137
$ws.Range("T4").Value = 0.007050342168100139

# Row 5
$ws.Range("G5").Value = 45.76217133333333
$ws.Range("H5").Value = 137.286514
$ws.Range("I5").Value = 0.6763939203605134
$ws.Range("J5").Value = 0.6763939203605135
$ws.Range("M5").Value = 2.105818666666667
$ws.Range("N5").Value = 6.317456
$ws.Range("O5").Value = 0.005641929445569353
$ws.Range("P5").Value = 0.005641929445569354
$ws.Range("Q5").Value = 96.36683462093156
$ws.Range("R5").Value = 867.3015115883841
$ws.Range("S5").Value = 0.003816166776086072
$ws.Range("T5").Value = 0.003816166776086073

# Row 6
$ws.Range("I6").Value = 0.1388778842960613
$ws.Range("J6").Value = 0.1388778842960613
$ws.Range("M6").Value = 361.2779286666667
$ws.Range("N6").Value = 1083.833786
$ws.Range("O6").Value = 0.9679392703861037
$ws.Range("P6").Value = 0.9679392703861038
$ws.Range("Q6").Value = 3394.544055112613
$ws.Range("R6").Value = 30550.89649601352
$ws.Range("S6").Value = 0.1344253579982953
$ws.Range("T6").Value = 0.1344253579982953

# Row 7
$ws.Range("I7").Value = 0.1388778842960613
$ws.Range("J7").Value = 0.1388778842960613
$ws.Range("O7").Value = 0.015995373883918
$ws.Range("P7").Value = 0.015995373883918
$ws.Range("S7").Value = 0.002221403683523004
$ws.Range("T7").Value = 0.002221403683523005

# Row 8
$ws.Range("I8").Value = 0.1388778842960613
$ws.Range("J8").Value = 0.1388778842960613
$ws.Range("M8").Value = 3.890485666666667
$ws.Range("N8").Value = 11.671457
$ws.Range("O8").Value = 0.01042342628440887
$ws.Range("P8").Value = 0.01042342628440887
$ws.Range("Q8").Value = 36.55475173926022
$ws.Range("R8").Value = 328.992765653342
$ws.Range("S8").Value = 0.00144758338949466
$ws.Range("T8").Value = 0.00144758338949466

# Row 9
$ws.Range("I9").Value = 0.1388778842960613
$ws.Range("J9").Value = 0.1388778842960613
$ws.Range("M9").Value = 2.105818666666667
$ws.Range("N9").Value = 6.317456
$ws.Range("O9").Value = 0.005641929445569353
$ws.Range("P9").Value = 0.005641929445569354
$ws.Range("Q9").Value = 19.78613601572622
$ws.Range("R9").Value = 178.075224141536
$ws.Range("S9").Value = 0.0007835392247483218
$ws.Range("T9").Value = 0.000783539224748322

# Row 10
$ws.Range("G10").Value = 12.29750866666667
$ws.Range("H10").Value = 36.892526
$ws.Range("I10").Value = 0.1817649787009828
$ws.Range("J10").Value = 0.1817649787009828
$ws.Range("M10").Value = 361.2779286666667
$ws.Range("N10").Value = 1083.833786
$ws.Range("O10").Value = 0.9679392703861037
$ws.Range("P10").Value = 0.9679392703861038
$ws.Range("Q10").Value = 4442.818458853716
$ws.Range("R10").Value = 39985.36612968345
$ws.Range("S10").Value = 0.1759374608655749
$ws.Range("T10").Value = 0.175937460865575

# Row 11
$ws.Range("G11").Value = 12.29750866666667
$ws.Range("H11").Value = 36.892526
$ws.Range("I11").Value = 0.1817649787009828
$ws.Range("J11").Value = 0.1817649787009828
$ws.Range("O11").Value = 0.015995373883918
$ws.Range("P11").Value = 0.015995373883918
$ws.Range("Q11").Value = 73.41838948159467
$ws.Range("R11").Value = 660.765505334352
$ws.Range("S11").Value = 0.002907398793324611
$ws.Range("T11").Value = 0.002907398793324611

# Row 12
$ws.Range("G12").Value = 12.29750866666667
$ws.Range("H12").Value = 36.892526
$ws.Range("I12").Value = 0.1817649787009828
$ws.Range("J12").Value = 0.1817649787009828
$ws.Range("M12").Value = 3.890485666666667
$ws.Range("N12").Value = 11.671457
$ws.Range("O12").Value = 0.01042342628440887
$ws.Range("P12").Value = 0.01042342628440887
$ws.Range("Q12").Value = 47.84328120337578
$ws.Range("R12").Value = 430.589530830382
$ws.Range("S12").Value = 0.001894613856576843
$ws.Range("T12").Value = 0.001894613856576843

# Row 13
$ws.Range("G13").Value = 12.29750866666667
$ws.Range("H13").Value = 36.892526
$ws.Range("I13").Value = 0.1817649787009828
$ws.Range("J13").Value = 0.1817649787009828
$ws.Range("M13").Value = 2.105818666666667
$ws.Range("N13").Value = 6.317456
$ws.Range("O13").Value = 0.005641929445569353
$ws.Range("P13").Value = 0.005641929445569354
$ws.Range("Q13").Value = 25.89632330376178
$ws.Range("R13").Value = 233.066909733856
$ws.Range("S13").Value = 0.001025505185506361
$ws.Range("T13").Value = 0.001025505185506361

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.2004796666666666
$ws.Range("H14").Value = 0.6014389999999999
$ws.Range("I14").Value = 0.002963216642442438
$ws.Range("J14").Value = 0.002963216642442439
$ws.Range("M14").Value = 361.2779286666667
$ws.Range("N14").Value = 1083.833786
$ws.Range("O14").Value = 0.9679392703861037
$ws.Range("P14").Value = 0.9679392703861038
$ws.Range("Q14").Value = 72.4288787131171
$ws.Range("R14").Value = 651.859908418054
$ws.Range("S14").Value = 0.002868213754881694
$ws.Range("T14").Value = 0.002868213754881694

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.2004796666666666
$ws.Range("H15").Value = 0.6014389999999999
$ws.Range("I15").Value = 0.002963216642442438
$ws.Range("J15").Value = 0.002963216642442439
$ws.Range("O15").Value = 0.015995373883918
$ws.Range("P15").Value = 0.015995373883918
$ws.Range("Q15").Value = 1.196900498258666
$ws.Range("R15").Value = 10.772104484328
$ws.Range("S15").Value = 0.00004739775809491495
$ws.Range("T15").Value = 0.00004739775809491496

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.2004796666666666
$ws.Range("H16").Value = 0.6014389999999999
$ws.Range("I16").Value = 0.002963216642442438
$ws.Range("J16").Value = 0.002963216642442439
$ws.Range("M16").Value = 3.890485666666667
$ws.Range("N16").Value = 11.671457
$ws.Range("O16").Value = 0.01042342628440887
$ws.Range("P16").Value = 0.01042342628440887
$ws.Range("Q16").Value = 0.7799632696247777
$ws.Range("R16").Value = 7.019669426623
$ws.Range("S16").Value = 0.00003088687023723231
$ws.Range("T16").Value = 0.00003088687023723233

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.2004796666666666
$ws.Range("H17").Value = 0.6014389999999999
$ws.Range("I17").Value = 0.002963216642442438
$ws.Range("J17").Value = 0.002963216642442439
$ws.Range("M17").Value = 2.105818666666667
$ws.Range("N17").Value = 6.317456
$ws.Range("O17").Value = 0.005641929445569353
$ws.Range("P17").Value = 0.005641929445569354
$ws.Range("Q17").Value = 0.4221738243537777
$ws.Range("R17").Value = 3.799564419184
$ws.Range("S17").Value = 0.00001671825922859714
$ws.Range("T17").Value = 0.00001671825922859715
